$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.484.64"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "1.871.41"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'316.57"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "'0.4666"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").Value = "'0.3759"
$ws.Range("E8").Value = "  +3.43%  "
$ws.Range("D9").Value = "'0.07407"
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("D10").Value = "'0.8922"
$ws.Range("E10").Value = "  +4.15%  "
$ws.Range("D11").Value = "'0.07972"
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("D12").Value = "'20.15"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "1.847.31"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "'5.458"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "'6.650"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "'93.03"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "'0.000008984"
$ws.Range("E18").Value = "  +4.96%  "
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("D21").Value = "27.509.88"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").Value = "'5.215"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").Value = "2.112.51"
$ws.Range("E24").Value = "  +7.99%  "
$ws.Range("D25").Value = "'153.10"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("D27").Value = "'18.64"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").Value = "'2.111"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").Value = "'5.222"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").Value = "'117.45"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "'0.08938"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "'0.7610"
$ws.Range("E32").Value = "  +6.68%  "
$ws.Range("D33").Value = "'2.995"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "'1.166"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").Value = "'4.519"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").Value = "'2.705"
$ws.Range("E36").Value = "  +12.26%  "
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01970"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05315"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").Value = "'3.004"
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("D41").Value = "'7.280"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").Value = "'0.5285"
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("D43").Value = "'0.1655"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("D44").Value = "'8.399"
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("D45").Value = "'0.4919"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("D46").Value = "'10.39"
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'104.29"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "'1.674"
$ws.Range("E49").Value = "  +3.91%  "
$ws.Range("D50").Value = "'0.06270"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'66.28"
$ws.Range("E51").Value = "  +3.86%  "
